$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    # Force the cell to remain plain text so numeric-looking strings
    # (e.g. "226.32", "0.557") are not silently converted to numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    # Reset to the default style so no stray formatting is introduced.
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "34.416.13"
Set-TextValue $ws.Range("E2") "  +0.66%  "
Set-TextValue $ws.Range("D3") "1.792.28"
Set-TextValue $ws.Range("E3") "  +0.21%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "226.32"
Set-TextValue $ws.Range("E5") "  -0.16%  "
Set-TextValue $ws.Range("D6") "0.557"
Set-TextValue $ws.Range("E6") "  +1.61%  "
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "32.63"
Set-TextValue $ws.Range("E8") "  +2.16%  "
Set-TextValue $ws.Range("D9") "0.296"
Set-TextValue $ws.Range("E9") "  +1.11%  "
Set-TextValue $ws.Range("D10") "0.0693"
Set-TextValue $ws.Range("E10") "  +0.51%  "
Set-TextValue $ws.Range("D11") "0.0949"
Set-TextValue $ws.Range("E11") "  +0.37%  "
Set-TextValue $ws.Range("D12") "2.052.56"
Set-TextValue $ws.Range("E12") "  +0.33%  "
Set-TextValue $ws.Range("D13") "1.795.54"
Set-TextValue $ws.Range("E13") "  +0.92%  "
Set-TextValue $ws.Range("D14") "11.03"
Set-TextValue $ws.Range("E14") "  +0.09%  "
Set-TextValue $ws.Range("D15") "0.635"
Set-TextValue $ws.Range("E15") "  +1.94%  "
Set-TextValue $ws.Range("D16") "34.422.59"
Set-TextValue $ws.Range("E16") "  +0.79%  "
Set-TextValue $ws.Range("E17") "  +2.16%  "
Set-TextValue $ws.Range("D18") "68.78"
Set-TextValue $ws.Range("E18") "  +0.69%  "
Set-TextValue $ws.Range("D19") "246.86"
Set-TextValue $ws.Range("E19") "  +0.21%  "
Set-TextValue $ws.Range("E20") "  +2.74%  "
Set-TextValue $ws.Range("E21") "  +3.65%  "
Set-TextValue $ws.Range("E22") "  -0.09%  "
Set-TextValue $ws.Range("D23") "4.15"
Set-TextValue $ws.Range("E23") "  +1.20%  "
Set-TextValue $ws.Range("E24") "  +1.36%  "
Set-TextValue $ws.Range("D25") "164.84"
Set-TextValue $ws.Range("E25") "  +2.29%  "
Set-TextValue $ws.Range("D26") "7.22"
Set-TextValue $ws.Range("E26") "  +0.71%  "
Set-TextValue $ws.Range("D27") "16.51"
Set-TextValue $ws.Range("E27") "  +0.98%  "
Set-TextValue $ws.Range("E28") "  +2.58%  "
Set-TextValue $ws.Range("E29") "  +0.08%  "
Set-TextValue $ws.Range("D30") "3.80"
Set-TextValue $ws.Range("E30") "  +3.61%  "
Set-TextValue $ws.Range("E31") "  +0.08%  "
Set-TextValue $ws.Range("D32") "0.0521"
Set-TextValue $ws.Range("E32") "  +0.49%  "
Set-TextValue $ws.Range("D33") "3.87"
Set-TextValue $ws.Range("E33") "  +6.53%  "
Set-TextValue $ws.Range("D34") "1.81"
Set-TextValue $ws.Range("E34") "  +0.83%  "
Set-TextValue $ws.Range("D35") "1.429.67"
Set-TextValue $ws.Range("E35") "  -0.92%  "
Set-TextValue $ws.Range("D36") "2.58"
Set-TextValue $ws.Range("E36") "  +7.11%  "
Set-TextValue $ws.Range("D37") "0.669"
Set-TextValue $ws.Range("E37") "  +3.12%  "
Set-TextValue $ws.Range("E38") "  +2.11%  "
Set-TextValue $ws.Range("E39") "  +0.07%  "
Set-TextValue $ws.Range("D40") "84.71"
Set-TextValue $ws.Range("E40") "  +5.17%  "
Set-TextValue $ws.Range("E41") "  +1.00%  "
Set-TextValue $ws.Range("D42") "0.936"
Set-TextValue $ws.Range("E42") "  +1.67%  "
Set-TextValue $ws.Range("E43") "  +1.83%  "
Set-TextValue $ws.Range("D44") "13.58"
Set-TextValue $ws.Range("E44") "  +0.62%  "
Set-TextValue $ws.Range("D45") "0.0523"
Set-TextValue $ws.Range("E46") "  +0.68%  "
Set-TextValue $ws.Range("E47") "  +0.09%  "
Set-TextValue $ws.Range("D48") "1.949.49"
Set-TextValue $ws.Range("E48") "  +0.09%  "
Set-TextValue $ws.Range("D49") "105.51"
Set-TextValue $ws.Range("E49") "  -0.31%  "
Set-TextValue $ws.Range("E50") "  -3.38%  "
Set-TextValue $ws.Range("E51") "  +0.00%  "

Write-Output "Updated cryptos list values on sheet '$($ws.Name)'"
